$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Fill in row 14 with a new time-tracking entry
$ws.Range("E14").Value = "LucaB"
$ws.Range("F14").Value = "GDPR"
$ws.Range("G14").Value = "Documentazione"
$ws.Range("H14").Value = (Get-Date -Year 2019 -Month 3 -Day 29 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("I14").Value = 60

# Update the active selection to match the edited cell
$ws.Range("I14").Select()
